$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 386, shifting existing rows 386:399 down to 387:400
$ws.Rows.Item(386).Insert()

# Populate the new row 386 with the new record's data
$ws.Range("A386").Value = 4
$ws.Range("B386").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C386").Value = "Los Lagos"
$ws.Range("D386").Value = 45075
$ws.Range("D386").NumberFormat = $ws.Range("D387").NumberFormat
$ws.Range("E386").Value = 10
$ws.Range("F386").Value = 100112044
$ws.Range("G386").Value = "Perejil"
$ws.Range("H386").Value = "Sin especificar"
$ws.Range("I386").Value = "Primera"
$ws.Range("J386").Value = 50
$ws.Range("K386").Value = 6000
$ws.Range("L386").Value = 6000
$ws.Range("M386").Value = 6000
$ws.Range("N386").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O386").Value = "Región de La Araucanía"
$ws.Range("P386").Value = 3000
$ws.Range("Q386").Value = 2
$ws.Range("R386").Value = "Hortaliza"
